# feat: add 2022-Q1 data
#
# - Insert a new sheet "2022-Q1" right before the "总计" (Total) sheet,
#   populated with the Q1-2022 per-fund holdings (cloned from the
#   "2021-Q4" sheet so it inherits the same column layout / styling).
# - Insert a new leading row into the "总计" summary sheet for the
#   2022-Q1 totals, pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q1" fund-holdings sheet.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")

# Copy "2021-Q4" to right before "总计" -- this clones headers, styles,
# dimensions etc. so the new sheet matches the existing per-quarter
# layout exactly.
$q4.Copy($total)

# Re-fetch "总计" -- its cached .Index doesn't update in place after the
# sheet collection changes, so look it up again to get the new position.
$total = $wb.Worksheets.Item("总计")
$q1_2022 = $wb.Worksheets.Item($total.Index - 1)
$q1_2022.Name = "2022-Q1"

# Force the numeric-looking text columns (fund code / fund size /
# position / weight / market value) to stay text, same as every other
# per-quarter sheet.
$q1_2022.Range("B2:G3").NumberFormat = "@"

# Row 2: 501007
$q1_2022.Range("A2").Value = 0
$q1_2022.Range("B2").Value = "501007"
$q1_2022.Range("C2").Value = "汇添富中证互联网医疗主题指数（LOF）A"
$q1_2022.Range("D2").Value = "0.58"
$q1_2022.Range("E2").Value = "93.89"
$q1_2022.Range("F2").Value = "6.67"
$q1_2022.Range("G2").Value = "0.0387"
$q1_2022.Range("H2").Value = 1

# Row 3: 501008
$q1_2022.Range("A3").Value = 1
$q1_2022.Range("B3").Value = "501008"
$q1_2022.Range("C3").Value = "汇添富中证互联网医疗主题指数（LOF）C"
$q1_2022.Range("D3").Value = "0.19"
$q1_2022.Range("E3").Value = "93.89"
$q1_2022.Range("F3").Value = "6.67"
$q1_2022.Range("G3").Value = "0.0127"
$q1_2022.Range("H3").Value = 1

# ---------------------------------------------------------------------
# 2) Insert the 2022-Q1 row at the top of the "总计" sheet's data.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows("2:2").Insert()
$total.Range("B2:D2").ClearFormats()

# A2 needs the same index-column style as the rows below it.
$total.Range("A3").Copy($total.Range("A2"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.05

# Renumber the 0-based index column (A) for every row now that a new
# row sits on top -- the row-insert shifts values down verbatim without
# renumbering them.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4

# Restore the original active sheet/selection (copying/inserting leaves
# the newly touched sheet focused).
$wb.Worksheets.Item("2021-Q1").Activate()
